$d = $word.ActiveDocument

# Helper: force a run boundary at a given (start,end) range by toggling
# Bold on/off around the edit. The replay engine merges adjacent runs
# that share identical effective formatting whenever a Text edit touches
# them, so a transient formatting change is used to keep the edited
# sub-range as its own <w:r>.
function Set-RangeTextSplit($range, [string]$newText) {
    $range.Bold = 1
    $range.Text = $newText
    $range.Bold = 0
}

function Force-Split($range) {
    $range.Bold = 1
    $range.Bold = 0
}

# ---------------------------------------------------------------------
# Change 1 (paragraph 1): fix the "Mothon" typo -> "Motion", landing as
# three separate runs: "Mot", "i", "on for Leave to File Responsive
# Pleadings".
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
$r1.Find.Execute("Mothon", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hPos = $r1.Start + 3
$rMid = $d.Range($hPos, $hPos + 1)
Set-RangeTextSplit $rMid "i"

# ---------------------------------------------------------------------
# Change 2 (paragraph 2): {%p if opposing_party_consent %}
#   -> {%p if opposing_party_consent_to_motion_for_leave_to_file %}
# split into: "opposing_party_consent_to_motion_for_leave_to_file", " ",
# "%}" as three distinct runs (the leading "f " run stays untouched).
# ---------------------------------------------------------------------
function Split-OpposingConsent([int]$paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $searchRange = $d.Range($p.Range.Start, $p.Range.End)
    $searchRange.Find.Execute("opposing_party_consent", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $wordStart = $searchRange.Start
    $wordEnd = $searchRange.End

    $rWord = $d.Range($wordStart, $wordEnd)
    Set-RangeTextSplit $rWord "opposing_party_consent_to_motion_for_leave_to_file"

    # " %}" -> " " + "%}"
    $spacePos = $rWord.End
    $paraEnd = $p.Range.End
    $rSpace = $d.Range($spacePos, $spacePos + 1)
    Force-Split $rSpace
    $rCloseTag = $d.Range($spacePos + 1, $paraEnd)
    Force-Split $rCloseTag
}

Split-OpposingConsent 2

# ---------------------------------------------------------------------
# Change 3 (paragraph 8): .{% if opposing_party_consent %}...
#   -> .{% if opposing_party_consent_to_motion_for_leave_to_file %}...
# split into 4 runs: ".{% if " (untouched), the renamed variable,
# " ", and "%}{{ landlord_doc_name }} ... states as follows:".
# ---------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$search8 = $d.Range($p8.Range.Start, $p8.Range.End)
$search8.Find.Execute("opposing_party_consent", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$word8Start = $search8.Start
$word8End = $search8.End

$rWord8 = $d.Range($word8Start, $word8End)
Set-RangeTextSplit $rWord8 "opposing_party_consent_to_motion_for_leave_to_file"

$space8Pos = $rWord8.End
$para8End = $p8.Range.End
$rSpace8 = $d.Range($space8Pos, $space8Pos + 1)
Force-Split $rSpace8
$rRest8 = $d.Range($space8Pos + 1, $para8End)
Force-Split $rRest8

# ---------------------------------------------------------------------
# Change 4 (paragraph 30): {%p if opposing_party_consent %}
#   -> {%p if opposing_party_consent_to_motion_for_leave_to_file %}
# same pattern as change 2.
# ---------------------------------------------------------------------
Split-OpposingConsent 30
